$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Find the last used row and the header row to locate the "Price" column.
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count
    $lastCol = $usedRange.Columns.Count

    # Header row is row 1; find the "Price" column index.
    $priceCol = 0
    for ($c = 1; $c -le $lastCol; $c++) {
        $headerValue = $ws.Cells.Item(1, $c).Value()
        if ($headerValue -eq "Price") {
            $priceCol = $c
        }
    }

    $totalRow = $lastRow + 1

    # Column A gets the label, column B gets the total of the Price column.
    $ws.Cells.Item($totalRow, 1).Value = "Monthly total: "

    $sum = 0.0
    for ($r = 2; $r -le $lastRow; $r++) {
        $cellValue = $ws.Cells.Item($r, $priceCol).Value()
        $sum = $sum + [double]$cellValue
    }
    $ws.Cells.Item($totalRow, 2).Value = $sum
}
